# Daily-Export-Template: insert a "Prediction" / "Daily Expected kWh" column
# between the existing "Weather Station" block (B:D) and the "Production"
# block (previously E:F, now shifted right to F:G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: merged sub-headers ------------------------------------------
# B6:D6 ("Weather Station") keeps its range/style untouched.

# The old "Production" merge (E6:F6) needs to move one column right (F6:G6);
# free it up first.
$ws.Range("E6:F6").UnMerge()

# E6 becomes the new standalone "Prediction" header (bold, centered,
# same look as B6/F6 but NOT merged).
$ws.Range("E6").Value2 = "Prediction"
$ws.Range("E6").Font.Bold = $true
$ws.Range("E6").HorizontalAlignment = -4108   # xlCenter

# F6:G6 becomes "Production" (previously E6:F6).
$ws.Range("F6").Value2 = "Production"
$ws.Range("G6").Value2 = ""
$ws.Range("F6:G6").Merge()
$ws.Range("F6:G6").Font.Bold = $true
$ws.Range("F6:G6").HorizontalAlignment = -4108   # xlCenter

# --- Row 7: column headers ----------------------------------------------
# Shift "Inverter Production kWh" / "Gen Meter Reading kWh" one column
# right; they keep their existing (wrap-only) look.
$ws.Range("G7").Value2 = $ws.Range("F7").Value2
$ws.Range("G7").Font.Bold = $true
$ws.Range("G7").WrapText = $true

$ws.Range("F7").Value2 = "Inverter Production`nkWh"
$ws.Range("F7").Font.Bold = $true
$ws.Range("F7").WrapText = $true

# New E7 header: "Daily Expected kWh" - wrapped AND centered.
$ws.Range("E7").Value2 = "Daily Expected`nkWh"
$ws.Range("E7").Font.Bold = $true
$ws.Range("E7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E7").WrapText = $true

# --- Column widths: new column G matches the others (B:F) ---------------
$ws.Columns("G:G").ColumnWidth = $ws.Columns("F:F").ColumnWidth

# --- Selection, matching the saved workbook state ------------------------
$ws.Range("E7").Select()
